$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new expense row (row 6)
$ws.Range("A6").Value = "digikey prototype order"
$ws.Range("B6").Value = "lots!"
$ws.Range("C6").Value = 50
$ws.Range("D6").Value = "digikey"

# Match the format used by other price cells in column C (currency, no decimals)
$ws.Range("C6").NumberFormat = $ws.Range("C3").NumberFormat

# Update selection to B9, as in the saved file
$ws.Range("B9").Select()

$wb.Save()
